# Update "想去人数" (wish-to-attend count, column F) figures on the
# 展览 (Exhibition) and 全部类型 (All Types) sheets, matching the refreshed
# scrape output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 40
$ws1.Range("F4").Value = 1422
$ws1.Range("F6").Value = 1040
$ws1.Range("F7").Value = 10784
$ws1.Range("F12").Value = 723
$ws1.Range("F13").Value = 12101
$ws1.Range("F14").Value = 12564

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 40
$ws4.Range("F5").Value = 1422
$ws4.Range("F7").Value = 1040
$ws4.Range("F8").Value = 10784
$ws4.Range("F13").Value = 723
$ws4.Range("F14").Value = 12101
$ws4.Range("F15").Value = 12564
